$p = $ppt.ActivePresentation

# -------------------------------------------------------------------------
# Slide 1: "Content Placeholder 9" - update group-member names.
# -------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 2: "Pratik Darade (A20536376)" -> "Ishan Prabhune"
$para2 = $tr.Paragraphs(2)
$para2.Text = "Ishan Prabhune"

# Paragraph 3: "Harshal Sawant (A20538827)" -> "Pratik " / "Darade" / " "
$para3 = $tr.Paragraphs(3)
$para3.Text = "Pratik "
[void]$para3.InsertAfter("Darade")
[void]$para3.InsertAfter(" ")

# Paragraph 4: "Ishan Prabhune (A20538828)" -> "Harshal " / "Sawant" / " "
$para4 = $tr.Paragraphs(4)
$para4.Text = "Harshal "
[void]$para4.InsertAfter("Sawant")
[void]$para4.InsertAfter(" ")

# New trailing (5th) paragraph containing a single space.
[void]$tr.InsertAfter([char]13 + " ")

# -------------------------------------------------------------------------
# Slide 11: tidy up the run fragmentation in the title and body text.
# -------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)

# Title shape: merge "and Satisfaction of " + "passenger's" into one run.
$titleShp = $s11.Shapes.Item(1)
$titleTr = $titleShp.TextFrame.TextRange
$f1 = $titleTr.Find("and Satisfaction of passenger")
$merge1 = $titleTr.Characters($f1.Start, $f1.Length + 2)
$merge1.Text = "and Satisfaction of passenger" + [char]8217 + "s"

# Body shape: merge the split runs back together.
$bodyShp = $s11.Shapes.Item(2)
$bodyTr = $bodyShp.TextFrame.TextRange

# "delay of Arrival (" + "Arrival Delay" + ")." -> one run
$f2 = $bodyTr.Find("delay of Arrival (")
$merge2 = $bodyTr.Characters($f2.Start, 33)
$merge2.Text = "delay of Arrival (Arrival Delay)."

# "level" + " " -> "level "
$f3 = $bodyTr.Find("level")
$merge3 = $bodyTr.Characters($f3.Start, 6)
$merge3.Text = "level "

# "Satisfaction (" + "Satisfaction" + ")." -> one run
$f4 = $bodyTr.Find("Satisfaction (")
$merge4 = $bodyTr.Characters($f4.Start, 28)
$merge4.Text = "Satisfaction (Satisfaction)."
